# DescripcionListado.xlsx — convert "Reactivos" catalog sheet into the
# generic "Catálogos" catalog sheet: rename sheet/defined name, swap the
# Contaq-system columns (Clave Contaq / Nombre Contaq) for a single
# "Descripción" column, and drop the old "Activo" column from E, moving
# "Activo" into column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet; Excel keeps the defined name's RefersTo formula
# pointed at the sheet automatically once we rename it below.
$ws.Name = "Catálogos"

# Rename + repoint the workbook-scoped defined name to the new sheet name
# and the shrunk A4:D5 range (column E is going away).
$dn = $wb.Names.Item(1)
$dn.Name = "Catalogos"
$dn.RefersTo = "=Catálogos!`$A`$4:`$D`$5"

# Header row (row 3): replace "Clave Contaq" / "Nombre Contaq" with a
# single "Descripción" header, and move "Activo" into column D.
$ws.Range("C3").Value = "Descripción"
$ws.Range("D3").Value = "Activo"

# Template row (row 4): same swap for the merge-field placeholders.
$ws.Range("C4").Value = "{{item.Descripcion}}"
$ws.Range("D4").Value = "{{item.Activo}}"

# Column E is no longer part of the table; fully clear both its format
# and its contents so the cells are dropped rather than left as empty
# styled placeholders.
$ws.Range("E3").ClearFormats()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()

# Leave the selection where the author left it when saving.
$null = $ws.Range("G1").Select()
